$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, shifting existing rows 3-7 down to 4-8
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly entry
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44490
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100101
$ws.Range("H3").Value = "Berries"
$ws.Range("I3").Value = 100101001
$ws.Range("J3").Value = "Arándano (blue)"
$ws.Range("K3").Value = "Sin especificar"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 400
$ws.Range("N3").Value = 9500
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 9750
$ws.Range("Q3").Value = "$/bandeja 2 kilos"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 4875
$ws.Range("T3").Value = 2
